$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualizei dados bibi: linha 8 (ano 2025) com os numeros revisados
$ws.Range("C8").Value = 1205
$ws.Range("D8").Value = 197
$ws.Range("E8").Value = 1008
$ws.Range("F8").Value = 8.08039376538146
$ws.Range("G8").Value = 83.65145228215768
$ws.Range("H8").Value = 16.34854771784232
